$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row (2-10) for columns E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$data = @{
    2  = @{ E=3; F=1; G=37.833119; H=113.499357; I=0.2771305381131279; J=0.2771305381131279; K=3; L=1; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=4428.160880281997; R=39853.44792253798; S=0.0899389927691269; T=0.0899389927691269 }
    3  = @{ E=3; F=1; G=37.833119; H=113.499357; I=0.2771305381131279; J=0.2771305381131279; K=3; L=1; M=101.5800373333333; N=304.740112; O=0.281657135515876; P=0.281657135515876; Q=3843.089640456443; R=34587.80676410798; S=0.07805579352891692; T=0.07805579352891689 }
    4  = @{ E=3; F=1; G=37.833119; H=113.499357; I=0.2771305381131279; J=0.2771305381131279; K=3; L=1; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=5373.316422035931; R=48359.84779832338; S=0.1091357518150841; T=0.1091357518150841 }
    5  = @{ E=3; F=1; G=63.88336466666667; H=191.650094; I=0.4679506129682439; J=0.467950612968244;  K=3; L=1; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=7477.200500379641; R=67294.80450341677; S=0.1518670843084027; T=0.1518670843084027 }
    6  = @{ E=3; F=1; G=63.88336466666667; H=191.650094; I=0.4679506129682439; J=0.467950612968244;  K=3; L=1; M=101.5800373333333; N=304.740112; O=0.281657135515876; P=0.281657135515876; Q=6489.274567818949; R=58403.47111037053; S=0.1318016292115339; T=0.1318016292115339 }
    7  = @{ E=3; F=1; G=63.88336466666667; H=191.650094; I=0.4679506129682439; J=0.467950612968244;  K=3; L=1; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=9073.149175417178; R=81658.3425787546; S=0.1842818994483074; T=0.1842818994483074 }
    8  = @{ E=3; F=1; G=34.80083866666666; H=104.402516; I=0.2549188489186281; J=0.2549188489186282; K=3; L=1; M=117.044563; N=351.133689; O=0.3245365645427815; P=0.3245365645427815; Q=4073.248953773502; R=36659.24058396152; S=0.08273048746525193; T=0.08273048746525194 }
    9  = @{ E=3; F=1; G=34.80083866666666; H=104.402516; I=0.2549188489186281; J=0.2549188489186282; K=3; L=1; M=101.5800373333333; N=304.740112; O=0.281657135515876; P=0.281657135515876; Q=3535.07049099131; R=31815.63441892179; S=0.07179971277542517; T=0.07179971277542517 }
    10 = @{ E=3; F=1; G=34.80083866666666; H=104.402516; I=0.2549188489186281; J=0.2549188489186282; K=3; L=1; M=142.0267893333333; N=426.080368; O=0.3938062999413425; P=0.3938062999413425; Q=4942.651381933987; R=44483.86243740589; S=0.1003886486779511; T=0.1003886486779511 }
}

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $cols) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
